$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text-like numeric string into a cell without Excel
# re-interpreting it as a Number (so the cell keeps type=Text, as the
# source file stores these as inline/shared strings) and without
# attaching a new (text) number-format style to the cell.
# Trick: put a formula that evaluates to the literal string, then
# Copy + PasteSpecial(values-only) it onto itself - this "bakes in"
# the text result and drops the formula, leaving a plain string cell
# with its original (unstyled) formatting untouched.
function Set-TextValue($range, $text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "29.304.53"
$ws.Range("E2").Value = "  -2.10%  "
Set-TextValue $ws.Range("D3") "1.855.80"
$ws.Range("E3").Value = "  -1.14%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "239.48"
$ws.Range("E5").Value = "  -1.18%  "
Set-TextValue $ws.Range("D6") "0.6953"
$ws.Range("E6").Value = "  -6.47%  "
$ws.Range("E7").Value = "  +0.13%  "
Set-TextValue $ws.Range("D8") "0.3071"
$ws.Range("E8").Value = "  -2.49%  "
Set-TextValue $ws.Range("D9") "0.07567"
$ws.Range("E9").Value = "  +4.78%  "
Set-TextValue $ws.Range("D10") "23.76"
$ws.Range("E10").Value = "  -3.76%  "
Set-TextValue $ws.Range("D11") "0.08112"
$ws.Range("E11").Value = "  -3.23%  "
Set-TextValue $ws.Range("D12") "1.856.89"
$ws.Range("E12").Value = "  -1.84%  "
Set-TextValue $ws.Range("D13") "0.7262"
$ws.Range("E13").Value = "  -3.55%  "
Set-TextValue $ws.Range("D14") "5.217"
$ws.Range("E14").Value = "  -3.80%  "
Set-TextValue $ws.Range("D15") "89.30"
$ws.Range("E15").Value = "  -3.56%  "
Set-TextValue $ws.Range("D16") "29.424.39"
$ws.Range("E16").Value = "  -1.69%  "
Set-TextValue $ws.Range("D17") "5.902"
$ws.Range("E17").Value = "  -2.98%  "
Set-TextValue $ws.Range("D18") "242.54"
$ws.Range("E18").Value = "  -4.34%  "
Set-TextValue $ws.Range("D19") "0.000007768"
$ws.Range("E19").Value = "  -1.12%  "
Set-TextValue $ws.Range("D20") "13.14"
$ws.Range("E20").Value = "  -3.40%  "
Set-TextValue $ws.Range("D22") "2.122.94"
$ws.Range("E22").Value = "  -0.23%  "
Set-TextValue $ws.Range("D23") "1.001"
$ws.Range("E23").Value = "  +0.05%  "
Set-TextValue $ws.Range("D24") "7.625"
$ws.Range("E24").Value = "  -5.05%  "
Set-TextValue $ws.Range("D25") "9.075"
$ws.Range("E25").Value = "  -2.11%  "
Set-TextValue $ws.Range("D26") "162.50"
$ws.Range("E26").Value = "  -1.58%  "
Set-TextValue $ws.Range("D27") "0.1466"
$ws.Range("E27").Value = "  -6.23%  "
Set-TextValue $ws.Range("D28") "18.11"
$ws.Range("E28").Value = "  -3.12%  "
Set-TextValue $ws.Range("D29") "1.935"
$ws.Range("E29").Value = "  -4.94%  "
Set-TextValue $ws.Range("D30") "1.404"
$ws.Range("E30").Value = "  -7.48%  "
Set-TextValue $ws.Range("D31") "1.515"
$ws.Range("E31").Value = "  -1.31%  "
Set-TextValue $ws.Range("D32") "4.427"
$ws.Range("E32").Value = "  -3.76%  "
Set-TextValue $ws.Range("D33") "4.053"
$ws.Range("E33").Value = "  -5.32%  "
Set-TextValue $ws.Range("D34") "0.05233"
$ws.Range("E34").Value = "  -1.81%  "
Set-TextValue $ws.Range("D35") "1.197"
$ws.Range("E35").Value = "  -3.18%  "
Set-TextValue $ws.Range("D36") "0.7182"
$ws.Range("E36").Value = "  -4.19%  "
Set-TextValue $ws.Range("D37") "1.001"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -0.89%  "
Set-TextValue $ws.Range("D39") "0.01867"
$ws.Range("E39").Value = "  -4.98%  "
Set-TextValue $ws.Range("D40") "2.709"
$ws.Range("E40").Value = "  -1.74%  "
Set-TextValue $ws.Range("D41") "0.8829"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("E42").Value = "  -4.93%  "
Set-TextValue $ws.Range("D43") "5.877"
$ws.Range("E43").Value = "  -3.29%  "
Set-TextValue $ws.Range("D44") "69.94"
$ws.Range("E44").Value = "  -3.49%  "
Set-TextValue $ws.Range("D45") "1.046.84"
$ws.Range("E45").Value = "  -5.94%  "
$ws.Range("E46").Value = "  +0.08%  "
Set-TextValue $ws.Range("D47") "102.91"
Set-TextValue $ws.Range("D48") "7.273"
$ws.Range("E48").Value = "  -4.55%  "
Set-TextValue $ws.Range("D49") "2.018.54"
$ws.Range("E49").Value = "  -0.33%  "
Set-TextValue $ws.Range("D50") "1.743"
$ws.Range("E50").Value = "  -6.06%  "
Set-TextValue $ws.Range("D51") "9.277"
$ws.Range("E51").Value = "  -2.15%  "
